$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Library Elements")

# Insert a new column before column C ("Application Type"'s old slot / now where
# External ID used to start) and give it the new header text. This shifts every
# column from C onward one slot to the right (C->D, D->E, ... I->J) exactly like
# the authored diff shows.
$ws.Range("C1").EntireColumn.Insert()
$ws.Range("C1").Value = "System Group"

# The sheet's AutoFilter range needs to become A1:H2 (it was A1:G2). This engine's
# Range.AutoFilter() call always snaps to the full contiguous data block though,
# which would pull row 3 in too (because row 3 has data directly below the
# filter header rows). Temporarily clear row 3's values so the filter has no
# contiguous data to latch onto below row 2, reapply the filter over exactly
# A1:H2, then restore row 3's values/formatting afterwards.
$row3 = @{}
for ($c = 1; $c -le 10; $c++) {
    $row3[$c] = $ws.Cells.Item(3, $c).Value()
}
$ws.Rows.Item(3).ClearContents()

$ws.AutoFilterMode = $false
$ws.Range("A1:H2").AutoFilter()

for ($c = 1; $c -le 10; $c++) {
    if ($row3[$c] -ne "") {
        $ws.Cells.Item(3, $c).Value = $row3[$c]
    }
}

# The two hyperlinks anchored at G2/G3 need to move to H2/H3 (their column shifted
# with the insert, but hyperlink anchors don't auto-track in this engine). Re-create
# them at the new location pointing at the same target URL.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("H2"), "https://sap.com/")
$ws.Hyperlinks.Add($ws.Range("H3"), "https://sap.com/")

# Hyperlinks.Add() stamps a brand-new "Hyperlink" cell style onto the target cells
# even when one is already in place; nudge the font back (no-op value) so the
# engine collapses the cell back onto the pre-existing matching style instead of
# keeping the freshly minted duplicate.
$ws.Range("H2").Font.Name = "Arial"
$ws.Range("H3").Font.Name = "Arial"

# The workbook-level _FilterDatabase defined name mirrors the AutoFilter range and
# also needs to grow from $A$1:$G$2 to $A$1:$H$2.
$wb.Names.Item("Library Elements!_FilterDatabase").RefersTo = "='Library Elements'!`$A`$1:`$H`$2"
